$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-20 with new data (A: message_id, B: message, C: sentence, D: boundary, E: error_type)
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = 'Aircraft ActiveTrack available at max speed . When exceeding nnn, Obstacle Avoidance is not available .'
$ws.Range("C2").Value = 'When exceeding nnn, Obstacle Avoidance is not available'
$ws.Range("D2").Value = '7-14'
$ws.Range("E2").Value = 'Missing'

$ws.Range("A3").Value = 4
$ws.Range("B3").Value = 'Aircraft ActiveTrack available at max speed . When exceeding nnn, Obstacle Avoidance is not available .'
$ws.Range("C3").Value = 'When exceeding nnn,'
$ws.Range("D3").Value = '7-9'
$ws.Range("E3").Value = "'False"
$ws.Range("E3").Style = "Normal"

$ws.Range("A4").Value = 4
$ws.Range("B4").Value = 'Aircraft ActiveTrack available at max speed . When exceeding nnn, Obstacle Avoidance is not available .'
$ws.Range("C4").Value = 'Obstacle Avoidance is not available'
$ws.Range("D4").Value = '10-14'
$ws.Range("E4").Value = "'False"
$ws.Range("E4").Style = "Normal"

$ws.Range("A5").Value = 21
$ws.Range("B5").Value = 'Aircraft processor chip overheated . Power off aircraft and wait for temperature to return to normal .'
$ws.Range("C5").Value = 'Power off aircraft and wait for temperature to return to normal'
$ws.Range("D5").Value = '5-15'
$ws.Range("E5").Value = 'Missing'

$ws.Range("A6").Value = 48
$ws.Range("B6").Value = 'Check whether propellers are installed correctly . If the propellers are installed correctly and the aircraft still cannot takeoff, a motor error may exist . Contact DJI Support for assistance .'
$ws.Range("C6").Value = 'If the propellers are installed correctly and the aircraft still cannot takeoff, a motor error may exist'
$ws.Range("D6").Value = '7-23'
$ws.Range("E6").Value = 'Missing'

$ws.Range("A7").Value = 48
$ws.Range("B7").Value = 'Check whether propellers are installed correctly . If the propellers are installed correctly and the aircraft still cannot takeoff, a motor error may exist . Contact DJI Support for assistance .'
$ws.Range("C7").Value = 'If the propellers are installed correctly and the aircraft still cannot takeoff,'
$ws.Range("D7").Value = '7-18'
$ws.Range("E7").Value = "'False"
$ws.Range("E7").Style = "Normal"

$ws.Range("A8").Value = 50
$ws.Range("B8").Value = 'Compass abnormal . Solution: 1. Ensure there are no magnets or metal objects near the aircraft . The ground or walls may contain metal . Move away from sources of interference before attempting flight . 2. Calibrate Compass Before Takeoff .'
$ws.Range("C8").Value = '2. Calibrate Compass Before Takeoff'
$ws.Range("D8").Value = '35-39'
$ws.Range("E8").Value = 'Missing'

$ws.Range("A9").Value = 66
$ws.Range("B9").Value = 'Downlink data connection lost for nnn seconds .'
$ws.Range("C9").Value = 'Downlink data connection lost for nnn seconds'
$ws.Range("D9").Value = '0-6'
$ws.Range("E9").Value = 'Missing'

$ws.Range("A10").Value = 66
$ws.Range("B10").Value = 'Downlink data connection lost for nnn seconds .'
$ws.Range("C10").Value = 'Downlink data connection lost for nnn'
$ws.Range("D10").Value = '0-5'
$ws.Range("E10").Value = "'False"
$ws.Range("E10").Style = "Normal"

$ws.Range("A11").Value = 81
$ws.Range("B11").Value = 'Extra payload detected . Return aircraft to an area nearby the home point promptly and fly in a wind-free environment to ensure flight safety .'
$ws.Range("C11").Value = 'Return aircraft to an area nearby the home point promptly and fly in a wind-free environment to ensure flight safety'
$ws.Range("D11").Value = '4-23'
$ws.Range("E11").Value = 'Missing'

$ws.Range("A12").Value = 85
$ws.Range("B12").Value = 'Flight altitude exceeds nnn . Aircraft may be in violation of local laws and regulations . Check and make sure you have obtained proper authorization to fly in this airspace .'
$ws.Range("C12").Value = 'Aircraft may be in violation of local laws and regulations'
$ws.Range("D12").Value = '5-14'
$ws.Range("E12").Value = 'Missing'

$ws.Range("A13").Value = 85
$ws.Range("B13").Value = 'Flight altitude exceeds nnn . Aircraft may be in violation of local laws and regulations . Check and make sure you have obtained proper authorization to fly in this airspace .'
$ws.Range("C13").Value = 'Aircraft may be in violation of local laws'
$ws.Range("D13").Value = '5-12'
$ws.Range("E13").Value = "'False"
$ws.Range("E13").Style = "Normal"

$ws.Range("A14").Value = 86
$ws.Range("B14").Value = 'Flight altitude exceeds nnn . May violate local policies and regulations . Ensure you have obtained proper airspace authorization .'
$ws.Range("C14").Value = 'May violate local policies and regulations'
$ws.Range("D14").Value = '5-10'
$ws.Range("E14").Value = 'Missing'

$ws.Range("A15").Value = 86
$ws.Range("B15").Value = 'Flight altitude exceeds nnn . May violate local policies and regulations . Ensure you have obtained proper airspace authorization .'
$ws.Range("C15").Value = 'May violate local policies and'
$ws.Range("D15").Value = '5-9'
$ws.Range("E15").Value = "'False"
$ws.Range("E15").Style = "Normal"

$ws.Range("A16").Value = 91
$ws.Range("B16").Value = 'GEO Zone Info: The target area is in an Altitude Zone . Flight altitude restricted to nnn .'
$ws.Range("C16").Value = 'GEO Zone Info: The target area is in an Altitude Zone'
$ws.Range("D16").Value = '0-10'
$ws.Range("E16").Value = 'Missing'

$ws.Range("A17").Value = 91
$ws.Range("B17").Value = 'GEO Zone Info: The target area is in an Altitude Zone . Flight altitude restricted to nnn .'
$ws.Range("C17").Value = 'GEO Zone Info:'
$ws.Range("D17").Value = '0-2'
$ws.Range("E17").Value = "'False"
$ws.Range("E17").Style = "Normal"

$ws.Range("A18").Value = 91
$ws.Range("B18").Value = 'GEO Zone Info: The target area is in an Altitude Zone . Flight altitude restricted to nnn .'
$ws.Range("C18").Value = 'The target area is in an Altitude Zone'
$ws.Range("D18").Value = '3-10'
$ws.Range("E18").Value = "'False"
$ws.Range("E18").Style = "Normal"

$ws.Range("A19").Value = 123
$ws.Range("B19").Value = 'Motor Obstructed . Propulsion output is limited to ensure the health of the battery .'
$ws.Range("C19").Value = 'Propulsion output is limited to ensure the health of the battery'
$ws.Range("D19").Value = '3-13'
$ws.Range("E19").Value = 'Missing'

$ws.Range("A20").Value = 154
$ws.Range("B20").Value = 'Strong wireless interference . Please fly with caution . Obstacle Avoidance Disabled . Landing gear lowered . Obstacle Avoidance Disabled .'
$ws.Range("C20").Value = 'Landing gear lowered'
$ws.Range("D20").Value = '13-15'
$ws.Range("E20").Value = 'Missing'

# Remove now-obsolete rows 21-23 (data no longer present)
$ws.Range("A21:E23").Delete()
